# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.517.22"
$ws.Range("D3").Value = "'1.673.12"
$ws.Range("E3").Value = "  +2.43%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'219.48"
$ws.Range("E5").Value = "  +2.35%  "
$ws.Range("D6").Value = "'0.526"
$ws.Range("E6").Value = "  +1.56%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "'29.51"
$ws.Range("E8").Value = "  +3.21%  "
$ws.Range("D9").Value = "'0.264"
$ws.Range("E9").Value = "  +2.04%  "
$ws.Range("D10").Value = "'0.0633"
$ws.Range("E10").Value = "  +3.95%  "
$ws.Range("D11").Value = "'0.0905"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "'1.916.59"
$ws.Range("E12").Value = "  +2.62%  "
$ws.Range("D13").Value = "'1.673.24"
$ws.Range("E13").Value = "  +2.50%  "
$ws.Range("D14").Value = "'0.614"
$ws.Range("E14").Value = "  +8.93%  "
$ws.Range("D15").Value = "'10.17"
$ws.Range("E15").Value = "  +9.73%  "
$ws.Range("D16").Value = "'3.95"
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").Value = "'30.521.51"
$ws.Range("E17").Value = "  +1.85%  "
$ws.Range("D18").Value = "'66.19"
$ws.Range("E18").Value = "  +3.24%  "
$ws.Range("D19").Value = "'242.98"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "'0.0₃0719"
$ws.Range("E20").Value = "  +2.58%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'4.25"
$ws.Range("E22").Value = "  +2.93%  "
$ws.Range("D23").Value = "'9.97"
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("D24").Value = "'2.15"
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("D25").Value = "'157.60"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "'15.84"
$ws.Range("E26").Value = "  +2.01%  "
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("D28").Value = "'6.65"
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").Value = "'0.0494"
$ws.Range("E30").Value = "  +1.78%  "
$ws.Range("D31").Value = "'1.14"
$ws.Range("E31").Value = "  +2.66%  "
$ws.Range("D32").Value = "'3.46"
$ws.Range("E32").Value = "  +2.54%  "
$ws.Range("D33").Value = "'1.501.90"
$ws.Range("E33").Value = "  +5.55%  "
$ws.Range("D34").Value = "'3.28"
$ws.Range("E34").Value = "  +3.41%  "
$ws.Range("D35").Value = "'1.75"
$ws.Range("E35").Value = "  +6.59%  "
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("D37").Value = "'83.60"
$ws.Range("E37").Value = "  +10.09%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0178"
$ws.Range("E38").Value = "  +4.89%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'0.592"
$ws.Range("E39").Value = "  +7.16%  "
$ws.Range("E40").Value = "  -3.85%  "
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").Value = "'0.836"
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.0497"
$ws.Range("E43").Value = "  +1.94%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'1.97"
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("D47").Value = "'5.55"
$ws.Range("E47").Value = "  +3.65%  "
$ws.Range("D48").Value = "'51.16"
$ws.Range("E48").Value = "  -3.33%  "
$ws.Range("D49").Value = "'1.808.60"
$ws.Range("E49").Value = "  +1.88%  "
$ws.Range("D50").Value = "'94.37"
$ws.Range("E50").Value = "  +5.37%  "
$ws.Range("E51").Value = "  -0.32%  "
